$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.04585075083376
$ws.Range("C2").Value = 7.360301762513083
$ws.Range("D2").Value = 6.375954411398858
$ws.Range("E2").Value = 11.22302240966283
$ws.Range("F2").Value = 33.16469232162403
$ws.Range("I2").Value = 26.03093852400587
$ws.Range("K2").Value = 11.69856910908127
$ws.Range("M2").Value = 15.07257873823612
$ws.Range("N2").Value = 20.84667234438739
$ws.Range("B3").Value = 10.80082906769146
$ws.Range("C3").Value = 7.1337695453687
$ws.Range("D3").Value = 6.380980532588324
$ws.Range("E3").Value = 11.00936252128347
$ws.Range("F3").Value = 33.03285283925451
$ws.Range("I3").Value = 26.03275416418846
$ws.Range("K3").Value = 11.52957008164124
$ws.Range("M3").Value = 14.91521557112774
$ws.Range("N3").Value = 20.89536727138223
$ws.Range("B4").Value = 10.6508510078446
$ws.Range("C4").Value = 6.99337644467318
$ws.Range("D4").Value = 6.384089450744963
$ws.Range("E4").Value = 10.87981425310023
$ws.Range("F4").Value = 32.95988283099728
$ws.Range("I4").Value = 26.03889924467055
$ws.Range("K4").Value = 11.42787674318543
$ws.Range("M4").Value = 14.82201611513312
$ws.Range("N4").Value = 20.92715163685552
$ws.Range("B5").Value = 10.58994961172879
$ws.Range("C5").Value = 6.93593923003781
$ws.Range("D5").Value = 6.385362302922467
$ws.Range("E5").Value = 10.82751168497949
$ws.Range("F5").Value = 32.93217201055971
$ws.Range("I5").Value = 26.04266711946303
$ws.Range("K5").Value = 11.38701150192759
$ws.Range("M5").Value = 14.7849382198413
$ws.Range("N5").Value = 20.94057831919072
$ws.Range("B6").Value = 10.57985294228748
$ws.Range("C6").Value = 6.926391260935576
$ws.Range("D6").Value = 6.385574023689678
$ws.Range("E6").Value = 10.81885869001642
$ws.Range("F6").Value = 32.92769344303233
$ws.Range("I6").Value = 26.04336905068523
$ws.Range("K6").Value = 11.38026219553656
$ws.Range("M6").Value = 14.77883708267225
$ws.Range("N6").Value = 20.94283646129556
$ws.Range("B7").Value = 10.65002866439168
$ws.Range("C7").Value = 6.992602602544991
$ws.Range("D7").Value = 6.384106592574411
$ws.Range("E7").Value = 10.87910679910263
$ws.Range("F7").Value = 32.95950089161918
$ws.Range("I7").Value = 26.03894494494335
$ws.Range("K7").Value = 11.42732321847626
$ws.Range("M7").Value = 14.82151236614238
$ws.Range("N7").Value = 20.92733079318954
$ws.Range("B8").Value = 10.96132898094359
$ws.Range("C8").Value = 7.282522458805022
$ws.Range("D8").Value = 6.377682791312762
$ws.Range("E8").Value = 11.14905838544888
$ws.Range("F8").Value = 33.11758968645442
$ws.Range("I8").Value = 26.03051996356058
$ws.Range("K8").Value = 11.63990096958499
$ws.Range("M8").Value = 15.01763489671001
$ws.Range("N8").Value = 20.86307112705745
$ws.Range("B9").Value = 11.57128493093085
$ws.Range("C9").Value = 7.836600268300938
$ws.Range("D9").Value = 6.365257980602211
$ws.Range("E9").Value = 11.68807213035444
$ws.Range("F9").Value = 33.4899879482679
$ws.Range("I9").Value = 26.05394463608827
$ws.Range("K9").Value = 12.07078432531327
$ws.Range("M9").Value = 15.42748970908397
$ws.Range("N9").Value = 20.75201045912923
$ws.Range("B10").Value = 12.01373217978528
$ws.Range("C10").Value = 8.229830199128971
$ws.Range("D10").Value = 6.356221315688016
$ws.Range("E10").Value = 12.08552788236742
$ws.Range("F10").Value = 33.80020555910999
$ws.Range("I10").Value = 26.09552222169447
$ws.Range("K10").Value = 12.39263563988688
$ws.Range("M10").Value = 15.74137279413274
$ws.Range("N10").Value = 20.67951803343176
$ws.Range("B11").Value = 12.21272592008287
$ws.Range("C11").Value = 8.40479135985982
$ws.Range("D11").Value = 6.352127407666846
$ws.Range("E11").Value = 12.26575487718634
$ws.Range("F11").Value = 33.94892998160833
$ws.Range("I11").Value = 26.11971907741272
$ws.Range("K11").Value = 12.53951656960197
$ws.Range("M11").Value = 15.88635823635129
$ws.Range("N11").Value = 20.64851411875031
$ws.Range("B12").Value = 12.2876713526541
$ws.Range("C12").Value = 8.470413656892427
$ws.Range("D12").Value = 6.350579371116269
$ws.Range("E12").Value = 12.33384801025665
$ws.Range("F12").Value = 34.00630835720654
$ws.Range("I12").Value = 26.12963952796254
$ws.Range("K12").Value = 12.59514833144275
$ws.Range("M12").Value = 15.94152749151847
$ws.Range("N12").Value = 20.63705743884101
$ws.Range("B13").Value = 12.27154995313976
$ws.Range("C13").Value = 8.456309822895031
$ws.Range("D13").Value = 6.350912672058931
$ws.Range("E13").Value = 12.3191909330003
$ws.Range("F13").Value = 33.99390432991146
$ws.Range("I13").Value = 26.12746932547475
$ws.Range("K13").Value = 12.58316742032952
$ws.Range("M13").Value = 15.92963474674585
$ws.Range("N13").Value = 20.6395122149423
$ws.Range("B14").Value = 12.21890037476039
$ws.Range("C14").Value = 8.410203169928394
$ws.Range("D14").Value = 6.352000006038497
$ws.Range("E14").Value = 12.27136044431339
$ws.Range("F14").Value = 33.95362949928131
$ws.Range("I14").Value = 26.12052006899682
$ws.Range("K14").Value = 12.54409352177284
$ws.Range("M14").Value = 15.89089198920226
$ws.Range("N14").Value = 20.64756588190286
$ws.Range("B15").Value = 12.18659537043933
$ws.Range("C15").Value = 8.381877352659194
$ws.Range("D15").Value = 6.352666315120602
$ws.Range("E15").Value = 12.24204064014898
$ws.Range("F15").Value = 33.92909696109533
$ws.Range("I15").Value = 26.11636204184821
$ws.Range("K15").Value = 12.52015943706404
$ws.Range("M15").Value = 15.86719415273357
$ws.Range("N15").Value = 20.65253595117035
$ws.Range("B16").Value = 12.00067453443105
$ws.Range("C16").Value = 8.218311214408059
$ws.Range("D16").Value = 6.356489185915043
$ws.Range("E16").Value = 12.07373178104735
$ws.Range("F16").Value = 33.79063624652464
$ws.Range("I16").Value = 26.09404704614782
$ws.Range("K16").Value = 12.38304121657603
$ws.Range("M16").Value = 15.73193764151432
$ws.Range("N16").Value = 20.6815839252436
$ws.Range("B17").Value = 11.88597494150126
$ws.Range("C17").Value = 8.116913904412762
$ws.Range("D17").Value = 6.358838586992823
$ws.Range("E17").Value = 11.97027961256535
$ws.Range("F17").Value = 33.70761941155265
$ws.Range("I17").Value = 26.0817093597056
$ws.Range("K17").Value = 12.29900324166696
$ws.Range("M17").Value = 15.64948984815147
$ws.Range("N17").Value = 20.69990936259456
$ws.Range("B18").Value = 11.81979292431543
$ws.Range("C18").Value = 8.058227384880032
$ws.Range("D18").Value = 6.360191504110616
$ws.Range("E18").Value = 11.91072638856341
$ws.Range("F18").Value = 33.66058829691498
$ws.Range("I18").Value = 26.07511047416123
$ws.Range("K18").Value = 12.25071426368463
$ws.Range("M18").Value = 15.60227797015983
$ws.Range("N18").Value = 20.71063539514303
$ws.Range("B19").Value = 11.79735136555719
$ws.Range("C19").Value = 8.038296549755765
$ws.Range("D19").Value = 6.360649860022829
$ws.Range("E19").Value = 11.89055630057124
$ws.Range("F19").Value = 33.64478868140622
$ws.Range("I19").Value = 26.07296169058817
$ws.Range("K19").Value = 12.23437433475747
$ws.Range("M19").Value = 15.58633040763115
$ws.Range("N19").Value = 20.71429894206964
$ws.Range("B20").Value = 11.89820723574975
$ws.Range("C20").Value = 8.127746159329737
$ws.Range("D20").Value = 6.358588324650218
$ws.Range("E20").Value = 11.98129800695613
$ws.Range("F20").Value = 33.71638262682157
$ws.Range("I20").Value = 26.08297125865576
$ws.Range("K20").Value = 12.30794473033943
$ws.Range("M20").Value = 15.65824519176667
$ws.Range("N20").Value = 20.69793936605092
$ws.Range("B21").Value = 12.23437656082073
$ws.Range("C21").Value = 8.423763462057936
$ws.Range("D21").Value = 6.3516805705132
$ws.Range("E21").Value = 12.28541417822948
$ws.Range("F21").Value = 33.96543071460164
$ws.Range("I21").Value = 26.12254069194118
$ws.Range("K21").Value = 12.55557061151372
$ws.Range("M21").Value = 15.9022648414732
$ws.Range("N21").Value = 20.6451926208969
$ws.Range("B22").Value = 12.45166364124569
$ws.Range("C22").Value = 8.613519207751539
$ws.Range("D22").Value = 6.347178906249209
$ws.Range("E22").Value = 12.48323786679785
$ws.Range("F22").Value = 34.13435745252306
$ws.Range("I22").Value = 26.15281621964338
$ws.Range("K22").Value = 12.71744844515398
$ws.Range("M22").Value = 16.06327574614072
$ws.Range("N22").Value = 20.61237400298446
$ws.Range("B23").Value = 12.33594033661317
$ws.Range("C23").Value = 8.512602981348014
$ws.Range("D23").Value = 6.349580408270064
$ws.Range("E23").Value = 12.37776391638437
$ws.Range("F23").Value = 34.04364620360487
$ws.Range("I23").Value = 26.13625451100997
$ws.Range("K23").Value = 12.63106571894828
$ws.Range("M23").Value = 15.97721758533531
$ws.Range("N23").Value = 20.6297385049454
$ws.Range("B24").Value = 11.89267775581724
$ws.Range("C24").Value = 8.122850112210982
$ws.Range("D24").Value = 6.358701461303365
$ws.Range("E24").Value = 11.97631682584284
$ws.Range("F24").Value = 33.7124186053504
$ws.Range("I24").Value = 26.08239921441778
$ws.Range("K24").Value = 12.30390219931338
$ws.Range("M24").Value = 15.65428631009755
$ws.Range("N24").Value = 20.69882940802636
$ws.Range("B25").Value = 11.40690118022264
$ws.Range("C25").Value = 7.688803054119468
$ws.Range("D25").Value = 6.368602145548027
$ws.Range("E25").Value = 11.54168738711924
$ws.Range("F25").Value = 33.38270164401219
$ws.Range("I25").Value = 26.04332712664073
$ws.Range("K25").Value = 11.95305621548011
$ws.Range("M25").Value = 15.31417686407625
$ws.Range("N25").Value = 20.78045583283098
